# Applies the "data up to 26" update to the survey-state workbook:
#   - fixes one previously-entered figure (G142)
#   - backfills the AR column for rows 120-124
#   - appends four new survey dates (rows 145-148) with their figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correction to an existing figure ---
$ws.Range("G142").Value = 0.3658553

# --- Backfill the AR column for rows 120-124 ---
$ws.Range("AR120").Value = 0.9615385
$ws.Range("AR121").Value = 0.9345793999999999
$ws.Range("AR122").Value = 0.9345793999999999
$ws.Range("AR123").Value = 0
$ws.Range("AR124").Value = 0.6289308

# --- Row 143 survey figures ---
$ws.Range("B143").Value = 0.4299176
$ws.Range("C143").Value = 0.5841029
$ws.Range("D143").Value = 0.7122694000000001
$ws.Range("F143").Value = 0.7868727
$ws.Range("G143").Value = 0.3841343
$ws.Range("H143").Value = 0.2769999
$ws.Range("I143").Value = 0.1875749
$ws.Range("J143").Value = 0.235953
$ws.Range("K143").Value = 0.4361181
$ws.Range("L143").Value = 0.5099305
$ws.Range("M143").Value = 0.51729
$ws.Range("O143").Value = 0.2928794
$ws.Range("P143").Value = 0.6668236
$ws.Range("Q143").Value = 0.5659569
$ws.Range("R143").Value = 0.3051378
$ws.Range("S143").Value = 0.58353
$ws.Range("T143").Value = 0.5800238
$ws.Range("U143").Value = 0.5266115
$ws.Range("V143").Value = 0.5770281
$ws.Range("W143").Value = 0.265331
$ws.Range("X143").Value = 0.2932868
$ws.Range("Y143").Value = 0.1983484
$ws.Range("Z143").Value = 0.2353787
$ws.Range("AA143").Value = 0.3484233
$ws.Range("AB143").Value = 0.3551902
$ws.Range("AD143").Value = 0.6270713
$ws.Range("AE143").Value = 0.4060952
$ws.Range("AF143").Value = 0.4098095
$ws.Range("AG143").Value = 0.2529333
$ws.Range("AH143").Value = 0.4274377
$ws.Range("AI143").Value = 0.2193734
$ws.Range("AJ143").Value = 0.3646814
$ws.Range("AK143").Value = 0.2534445
$ws.Range("AL143").Value = 0.5472644
$ws.Range("AM143").Value = 0.3251949
$ws.Range("AN143").Value = 0.3087605
$ws.Range("AO143").Value = 0.5384904
$ws.Range("AP143").Value = 0.3278763
$ws.Range("AQ143").Value = 0.2948958
$ws.Range("AS143").Value = 0.3278669
$ws.Range("AT143").Value = 0.6595007000000001
$ws.Range("AU143").Value = 0.4605847
$ws.Range("AV143").Value = 0.5583333
$ws.Range("AW143").Value = 0.7355757000000001
$ws.Range("AX143").Value = 0.5034583
$ws.Range("AY143").Value = 0.3731443
$ws.Range("BA143").Value = 0.40764
$ws.Range("BB143").Value = 0.3053555
$ws.Range("BC143").Value = 0.3369779
$ws.Range("BD143").Value = 0.3224896
$ws.Range("BE143").Value = 0.7194416

# --- Row 144 survey figures ---
$ws.Range("B144").Value = 0.5409677000000001
$ws.Range("C144").Value = 0.6563312
$ws.Range("D144").Value = 0.6758653
$ws.Range("F144").Value = 0.8682566
$ws.Range("G144").Value = 0.3732982
$ws.Range("H144").Value = 0.2715601
$ws.Range("I144").Value = 0.1603008
$ws.Range("J144").Value = 0.2424986
$ws.Range("K144").Value = 0.5102501
$ws.Range("L144").Value = 0.5558275
$ws.Range("M144").Value = 0.5254165
$ws.Range("O144").Value = 0.2913085
$ws.Range("P144").Value = 0.6317675
$ws.Range("Q144").Value = 0.5898221
$ws.Range("R144").Value = 0.3236772
$ws.Range("S144").Value = 0.5556495
$ws.Range("T144").Value = 0.5585957
$ws.Range("U144").Value = 0.5602354000000001
$ws.Range("V144").Value = 0.6214258
$ws.Range("W144").Value = 0.287014
$ws.Range("X144").Value = 0.3139187
$ws.Range("Y144").Value = 0.2650219
$ws.Range("Z144").Value = 0.2204673
$ws.Range("AA144").Value = 0.3137485
$ws.Range("AB144").Value = 0.3946895
$ws.Range("AD144").Value = 0.6780471
$ws.Range("AE144").Value = 0.442824
$ws.Range("AF144").Value = 0.4228497
$ws.Range("AG144").Value = 0.400355
$ws.Range("AH144").Value = 0.4924516
$ws.Range("AI144").Value = 0.256426
$ws.Range("AJ144").Value = 0.3672354
$ws.Range("AK144").Value = 0.2207373
$ws.Range("AL144").Value = 0.5672401
$ws.Range("AM144").Value = 0.3252889
$ws.Range("AN144").Value = 0.3048788
$ws.Range("AO144").Value = 0.5171485
$ws.Range("AP144").Value = 0.3007372
$ws.Range("AQ144").Value = 0.2739843
$ws.Range("AS144").Value = 0.3575959
$ws.Range("AT144").Value = 0.696344
$ws.Range("AU144").Value = 0.5482522
$ws.Range("AV144").Value = 0.506008
$ws.Range("AW144").Value = 0.7636775
$ws.Range("AX144").Value = 0.5754314
$ws.Range("AY144").Value = 0.3636353
$ws.Range("BA144").Value = 0.4546538
$ws.Range("BB144").Value = 0.3337814
$ws.Range("BC144").Value = 0.3328386
$ws.Range("BD144").Value = 0.3310033
$ws.Range("BE144").Value = 0.7952616

# --- Row 145 survey figures ---
$ws.Range("A145").Value = "23 06 2020"
$ws.Range("B145").Value = 0.3361238
$ws.Range("C145").Value = 0.6824983
$ws.Range("D145").Value = 0.6869927
$ws.Range("F145").Value = 0.8755231999999999
$ws.Range("G145").Value = 0.397048
$ws.Range("H145").Value = 0.3502029
$ws.Range("I145").Value = 0.1533982
$ws.Range("J145").Value = 0.244436
$ws.Range("K145").Value = 0.5442462
$ws.Range("L145").Value = 0.5817508
$ws.Range("M145").Value = 0.5425352
$ws.Range("O145").Value = 0.2230333
$ws.Range("P145").Value = 0.6353315
$ws.Range("Q145").Value = 0.5309932000000001
$ws.Range("R145").Value = 0.3323737
$ws.Range("S145").Value = 0.4743293
$ws.Range("T145").Value = 0.5715398
$ws.Range("U145").Value = 0.5476126
$ws.Range("V145").Value = 0.685158
$ws.Range("W145").Value = 0.317587
$ws.Range("X145").Value = 0.3285764
$ws.Range("Y145").Value = 0.2841275
$ws.Range("Z145").Value = 0.2611591
$ws.Range("AA145").Value = 0.335791
$ws.Range("AB145").Value = 0.4293275
$ws.Range("AD145").Value = 0.734653
$ws.Range("AE145").Value = 0.4720346
$ws.Range("AF145").Value = 0.4087258
$ws.Range("AG145").Value = 0.5940805
$ws.Range("AH145").Value = 0.4688323
$ws.Range("AI145").Value = 0.304383
$ws.Range("AJ145").Value = 0.3765047
$ws.Range("AK145").Value = 0.2634061
$ws.Range("AL145").Value = 0.5519790999999999
$ws.Range("AM145").Value = 0.3276885
$ws.Range("AN145").Value = 0.3154596
$ws.Range("AO145").Value = 0.5740784
$ws.Range("AP145").Value = 0.282927
$ws.Range("AQ145").Value = 0.3029777
$ws.Range("AS145").Value = 0.3449371
$ws.Range("AT145").Value = 0.7278261
$ws.Range("AU145").Value = 0.584093
$ws.Range("AV145").Value = 0.4992202
$ws.Range("AW145").Value = 0.8534987000000001
$ws.Range("AX145").Value = 0.5350799000000001
$ws.Range("AY145").Value = 0.3627693
$ws.Range("BA145").Value = 0.3484736
$ws.Range("BB145").Value = 0.297706
$ws.Range("BC145").Value = 0.3485606
$ws.Range("BD145").Value = 0.3954634
$ws.Range("BE145").Value = 0.6476771

# --- Row 146 survey figures ---
$ws.Range("A146").Value = "24 06 2020"
$ws.Range("B146").Value = 0.3527689
$ws.Range("C146").Value = 0.6726328
$ws.Range("D146").Value = 0.6582851
$ws.Range("F146").Value = 0.9453054
$ws.Range("G146").Value = 0.4242899
$ws.Range("H146").Value = 0.3440139
$ws.Range("I146").Value = 0.1581242
$ws.Range("J146").Value = 0.2820608
$ws.Range("K146").Value = 0.5092209
$ws.Range("L146").Value = 0.6522483
$ws.Range("M146").Value = 0.5396107999999999
$ws.Range("O146").Value = 0.2472188
$ws.Range("P146").Value = 0.5483575000000001
$ws.Range("Q146").Value = 0.5678628999999999
$ws.Range("R146").Value = 0.3402836
$ws.Range("S146").Value = 0.4992076
$ws.Range("T146").Value = 0.6622741
$ws.Range("U146").Value = 0.5391899
$ws.Range("V146").Value = 0.6737465
$ws.Range("W146").Value = 0.3009859
$ws.Range("X146").Value = 0.3503361
$ws.Range("Y146").Value = 0.2275811
$ws.Range("Z146").Value = 0.2584056
$ws.Range("AA146").Value = 0.3438407
$ws.Range("AB146").Value = 0.4526778
$ws.Range("AD146").Value = 0.8248751
$ws.Range("AE146").Value = 0.566879
$ws.Range("AF146").Value = 0.4120582
$ws.Range("AG146").Value = 0.8287523
$ws.Range("AH146").Value = 0.587907
$ws.Range("AI146").Value = 0.2964371
$ws.Range("AJ146").Value = 0.3599
$ws.Range("AK146").Value = 0.2570161
$ws.Range("AL146").Value = 0.5283232
$ws.Range("AM146").Value = 0.3327101
$ws.Range("AN146").Value = 0.2983306
$ws.Range("AO146").Value = 0.6271786
$ws.Range("AP146").Value = 0.3178141
$ws.Range("AQ146").Value = 0.277503
$ws.Range("AS146").Value = 0.344979
$ws.Range("AT146").Value = 0.6879987
$ws.Range("AU146").Value = 0.7117162
$ws.Range("AV146").Value = 0.5182799
$ws.Range("AW146").Value = 0.9088826
$ws.Range("AX146").Value = 0.5641392
$ws.Range("AY146").Value = 0.3315785
$ws.Range("BA146").Value = 0.286497
$ws.Range("BB146").Value = 0.2655729
$ws.Range("BC146").Value = 0.3601689
$ws.Range("BD146").Value = 0.4709291
$ws.Range("BE146").Value = 0.6664116

# --- Row 147 survey figures ---
$ws.Range("A147").Value = "25 06 2020"
$ws.Range("B147").Value = 0.3556092
$ws.Range("C147").Value = 0.6356416
$ws.Range("D147").Value = 0.6601846
$ws.Range("F147").Value = 1.0347193
$ws.Range("G147").Value = 0.4227923
$ws.Range("H147").Value = 0.3387272
$ws.Range("I147").Value = 0.1505218
$ws.Range("J147").Value = 0.2752599
$ws.Range("K147").Value = 0.499038
$ws.Range("L147").Value = 0.6783066
$ws.Range("M147").Value = 0.5795326
$ws.Range("O147").Value = 0.2374947
$ws.Range("P147").Value = 0.4828779
$ws.Range("Q147").Value = 0.6138405
$ws.Range("R147").Value = 0.3397252
$ws.Range("S147").Value = 0.4960112
$ws.Range("T147").Value = 0.5551384
$ws.Range("U147").Value = 0.4425484
$ws.Range("V147").Value = 0.724025
$ws.Range("W147").Value = 0.2804375
$ws.Range("X147").Value = 0.3708378
$ws.Range("Y147").Value = 0.2782856
$ws.Range("Z147").Value = 0.2201957
$ws.Range("AA147").Value = 0.34641
$ws.Range("AB147").Value = 0.5682218999999999
$ws.Range("AD147").Value = 0.9288629
$ws.Range("AE147").Value = 0.6107671
$ws.Range("AF147").Value = 0.4377089
$ws.Range("AG147").Value = 0.8198613
$ws.Range("AH147").Value = 0.523158
$ws.Range("AI147").Value = 0.3966866
$ws.Range("AJ147").Value = 0.3670605
$ws.Range("AK147").Value = 0.2897309
$ws.Range("AL147").Value = 0.5100969
$ws.Range("AM147").Value = 0.3565261
$ws.Range("AN147").Value = 0.3144984
$ws.Range("AO147").Value = 0.5892324
$ws.Range("AP147").Value = 0.3057701
$ws.Range("AQ147").Value = 0.2961598
$ws.Range("AS147").Value = 0.3243032
$ws.Range("AT147").Value = 0.7662813000000001
$ws.Range("AU147").Value = 0.7958097
$ws.Range("AV147").Value = 0.5864935999999999
$ws.Range("AW147").Value = 0.9468583
$ws.Range("AX147").Value = 0.6248169
$ws.Range("AY147").Value = 0.2945737
$ws.Range("BA147").Value = 0.1518602
$ws.Range("BB147").Value = 0.2638247
$ws.Range("BC147").Value = 0.3720608
$ws.Range("BD147").Value = 0.4922279
$ws.Range("BE147").Value = 0.7418856

# --- Row 148: date only (figures not yet available) ---
$ws.Range("A148").Value = "26 06 2020"

